$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the existing hyperlink on F2 before shifting rows down -------
# (Rows.Insert does not re-target existing Hyperlink objects to their new
# row, so leaving it in place would produce a stale hyperlink pointing at
# the now-empty F2 cell.)
$ws.Range("F2").Hyperlinks.Delete()

# --- Insert a new row above the existing data row (old row 2 -> row 3) ---
$ws.Rows.Item(2).Insert()

# --- New row 2: newly scraped listing -------------------------------------
$ws.Range("A2").Value = "2025-12-15 02:00:19"
$ws.Range("B2").Value = "【AIシステム構築】次のテストに向けた宿題自動出題システム"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("G2").Value = 318
$ws.Range("H2").Value = "🔥AI,Ai"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5453785")

# --- Row 3 (previously row 2): refresh timestamp and re-link its URL -----
$ws.Range("A3").Value = "2025-12-15 02:00:19"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5453611")

# --- Row 4: appended listing -----------------------------------------------
$ws.Range("A4").Value = "2025-12-15 02:00:19"
$ws.Range("B4").Value = "【急募】ホームページとLPの改善をお手伝いします!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("G4").Value = 18
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5453763")

# --- Widen column D by one character (27 -> 28) ---------------------------
# ColumnWidth is expressed in "characters of the Normal font" and is offset
# from the raw OOXML <col width> by 5/6 (0.8333...), so subtract that back
# out to land exactly on 28.
$ws.Columns.Item(4).ColumnWidth = 28 - 0.8333333333333334
